$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data values: Indice, Distancia, max, min, Tempo
$data = @(
    @(0, 1266, 1266, 1266, 0.01084530353546143),
    @(1, 911, 911, 911, 0.01063434282938639),
    @(2, 1233, 1233, 1233, 0.0124152421951294),
    @(3, 1430, 1430, 1430, 0.0123185396194458),
    @(4, 1198, 1198, 1198, 0.0126118262608846),
    @(5, 1284, 1284, 1284, 0.01239616870880127),
    @(6, 1545, 1545, 1545, 0.0126328706741333),
    @(7, 1733, 1733, 1733, 0.01213558514912923),
    @(8, 1770, 1770, 1770, 0.01232256094614665),
    @(9, 1741, 1741, 1741, 0.01234435240427653)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
